# ToDo.docx update
#
# 1. Remove the "Add in exceptions in code" bullet entirely.
# 2. Append a new bullet at the end of the list:
#      "Night School stop: should not set 'pathChosen' as a specific
#       number: change this"
#    (with the same list formatting as the preceding "Fix up ..." bullet,
#    "pathChosen" flagged the same way "LifeDetails" was flagged for
#    spell-check, and the _GoBack bookmark moved onto this new, final
#    paragraph.)

$d = $word.ActiveDocument

# --- Step 1: find & delete the "Add in exceptions in code" paragraph ---
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Add in exceptions in code*") {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -gt 0) {
    $d.Paragraphs.Item($targetIndex).Range.Delete()
}

# --- Step 2: drop the existing _GoBack bookmark (currently sitting on ---
# --- the "Fix up 'LifeDetails' Excel sheet" paragraph, now the last)  ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 3: append the new paragraph (with proofErr + bookmark) right ---
# --- after the current last paragraph, via a WordprocessingML fragment ---
$d2 = $word.ActiveDocument
$lastPara = $d2.Paragraphs.Item($d2.Paragraphs.Count)
$insertPoint = $d2.Range($lastPara.Range.End, $lastPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr>
<w:ilvl w:val="0"/>
<w:numId w:val="4"/>
</w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Night School stop: should not set &#8216;</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>pathChosen</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>&#8217; as a specific number: change this</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$insertPoint.InsertXML($xml)
